$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.822.87"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "2.810.10"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.42"
$ws.Range("E5").Value = "  +2.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.04"
$ws.Range("E6").Value = "  -4.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  +1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.33"
$ws.Range("E10").Value = "  -6.48%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.84"
$ws.Range("E13").Value = "  -1.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.72"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").Value = "3.246.70"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "2.804.27"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.927"
$ws.Range("E17").Value = "  +3.68%  "

$ws.Range("D18").Value = "51.535.31"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  +5.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.13"
$ws.Range("E20").Value = "  -2.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.33"
$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("D22").Value = "0.0₃0988"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.39"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.82"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0483"
$ws.Range("E30").Value = "  +16.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.50"
$ws.Range("E32").Value = "  +4.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.37"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.92"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("E35").Value = "  +11.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0845"
$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.23"
$ws.Range("E40").Value = "  -4.27%  "

$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "126.86"
$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.01"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  -2.13%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  -8.53%  "

$ws.Range("D46").Value = "2.073.31"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.86"
$ws.Range("E49").Value = "  +5.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.962"
$ws.Range("E50").Value = "  +1.62%  "
